$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data to the worksheet.
# NumberFormat is forced to text ("@") before assignment so that values
# such as "27.10" or "246.80" are not coerced into numeric types and
# lose their trailing zeroes / original string formatting.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.843.00"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +4.97%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.615.09"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.23%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.49%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.96"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.16%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.518"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +7.59%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.994"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.59%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "27.10"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +13.10%  "

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.97%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0601"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.08%  "

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.08%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.846.22"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.24%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.619.65"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.53%  "

# Row 14
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.543"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +6.76%  "

# Row 15
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.830.47"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +5.00%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.78"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.28%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "246.80"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +8.14%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.79"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +4.68%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.65"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.71%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0699"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.97%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.995"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.47%  "

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +4.84%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.31"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +4.66%  "

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +5.20%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.28"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.17%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.43"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +4.85%  "

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +6.16%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.44"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.73%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.995"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.47%  "

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.81%  "

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.96%  "

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.60%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.446.66"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.60%  "

# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.35%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.69%  "

# Row 36
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.52"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.63%  "

# Row 37
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.85"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +10.87%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.31"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.37%  "

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.50%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "56.66"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +32.15%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.538"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +5.80%  "

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.17%  "

# Row 43
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.802"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.17%  "

# Row 44
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.995"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.43%  "

# Row 45
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "67.76"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +9.82%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0469"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.57%  "

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.84%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.756.51"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.42%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.96"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.11%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.837"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.01%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0104"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.99%  "
